# Add files via upload
#
# "All structural tables have been made uniform in terms of how valency
# classes, X and Y columns etc. are organized. They still don't have
# interlinearized examples though."
#
# Column I is "X" and column J is "Y" in the valency table. A number of rows
# used the shorthand marker "TR" in the X column to mean "some nominative/
# accusative pair" and simply left the Y column blank instead of spelling the
# pattern out. This pass normalizes those rows: the X column becomes the
# explicit case "NOM" and the previously-blank Y column is filled in with the
# explicit case "ACC". Rows that used the other shorthand marker "*" (meaning
# "valency info not applicable / TBD") get that same "*" copied into the
# blank Y column as well, leaving X untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$colX = 9   # column I - "X"
$colY = 10  # column J - "Y"

$lastRow = $ws.Cells.Item(1, 1).SpecialCells(11).Row

for ($r = 2; $r -le $lastRow; $r++) {
    $xCell = $ws.Cells.Item($r, $colX)
    $yCell = $ws.Cells.Item($r, $colY)

    $xVal = $xCell.Value2
    $yVal = $yCell.Value2

    if ($yVal -eq $null -or $yVal -eq "") {
        if ($xVal -eq "TR") {
            $xCell.Value = "NOM"
            $yCell.Value = "ACC"
        }
        elseif ($xVal -eq "*") {
            $yCell.Value = "*"
        }
    }
}

# Update the view state left over from the author's editing session:
# scrolled so column E is the leftmost visible column, with M7 selected.
$ws.Range("E1").Select()
$ws.Range("M7").Select()
